$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date in column C for existing rows 2-28
# from 45537 (2024-09-02) to 45538 (2024-09-03).
for ($r = 2; $r -le 28; $r++) {
    $ws.Cells.Item($r, 3).Value = 45538
}

# Row 28 gains an explicit row height (15, custom height) in the new file.
$ws.Rows.Item(28).RowHeight = 15

# Append a new record in row 29.
$ws.Cells.Item(29, 1).Value = "A 36712-2024"
$ws.Cells.Item(29, 2).Value = 45537
$ws.Cells.Item(29, 2).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(29, 3).Value = 45538
$ws.Cells.Item(29, 3).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(29, 4).Value = "OKÄNT"
$ws.Cells.Item(29, 5).Value = "OKÄNT"
$ws.Cells.Item(29, 6).Value = "Övriga Aktiebolag"
$ws.Cells.Item(29, 7).Value = 6.4
$ws.Cells.Item(29, 8).Value = 0
$ws.Cells.Item(29, 9).Value = 0
$ws.Cells.Item(29, 10).Value = 0
$ws.Cells.Item(29, 11).Value = 0
$ws.Cells.Item(29, 12).Value = 0
$ws.Cells.Item(29, 13).Value = 0
$ws.Cells.Item(29, 14).Value = 0
$ws.Cells.Item(29, 15).Value = 0
$ws.Cells.Item(29, 16).Value = 0
$ws.Cells.Item(29, 17).Value = 0
$ws.Cells.Item(29, 18).WrapText = $true

Write-Output "done"
